$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("indicators")
$ws.Range("C2").Value = 'municipal administration'
$ws.Range("E2").Value = 'Sustainable plaid health goth pinterest YOLO, authentic hoodie hashtag fixie swag hella aesthetic banh mi fap fingerstache. Cred shoreditch godard, neutra deep v narwhal food truck flannel normcore.'
$ws.Range("C3").Value = 'municipal assembly'
$ws.Range("E3").Value = 'Waistcoat neutra synth 90''s distillery, +1 sartorial occupy pinterest ugh kickstarter.'
$ws.Range("C4").Value = 'cemëtery maintenancë'
$ws.Range("E4").Value = 'Kombucha sartorial try-hard offal food truck. Pitchfork YOLO brooklyn stumptown locavore.'
$ws.Range("C5").Value = 'çultural activities'
$ws.Range("E5").Value = 'Keytar swag portland, meh food truck quinoa knausgaard gastropub synth williamsburg. '
$ws.Range("C6").Value = 'protection of cultural heritage'
$ws.Range("E6").Value = 'Chia disrupt raw denim twee synth. Messenger bag waistcoat single-origin coffee, sartorial wolf vegan food truck. Hoodie narwhal freegan pop-up hammock hella wolf.'
$ws.Range("C7").Value = 'electricity supply'
$ws.Range("E7").Value = 'Narwhal XOXO cardigan, seitan microdosing pabst selvage cray lo-fi skateboard brunch. '
$ws.Range("C8").Value = 'emergency services'
$ws.Range("E8").Value = 'Gastropub art party sustainable, tattooed small batch vice leggings actually tote bag. Cray kinfolk hoodie mumblecore tumblr.'
$ws.Range("C9").Value = 'environmental protection'
$ws.Range("E9").Value = 'Cliche twee neutra bushwick, celiac photo booth distillery pitchfork selvage beard microdosing swag hoodie pinterest four dollar toast.'
$ws.Range("C10").Value = 'family medical centers'
$ws.Range("E10").Value = 'Ramps scenester shabby chic waistcoat food truck, four loko art party authentic.'
$ws.Range("C11").Value = 'management of municipal funds'
$ws.Range("E11").Value = 'Venmo celiac kale chips, distillery stumptown offal chillwave PBR&B twee cray ethical tote bag food truck. Tumblr DIY pour-over kinfolk vinyl. '
$ws.Range("C12").Value = 'public health'
$ws.Range("E12").Value = 'Fingerstache forage craft beer, intelligentsia chartreuse narwhal seitan fixie health goth paleo 3 wolf moon selvage single-origin coffee pinterest typewriter. '
$ws.Range("C13").Value = 'hospitals'
$ws.Range("E13").Value = '3 wolf moon stumptown bespoke, chillwave viral crucifix taxidermy tilde. Dreamcatcher deep v authentic, tote bag cliche shoreditch typewriter flannel disrupt.'
$ws.Range("C14").Value = 'public lighting'
$ws.Range("E14").Value = 'Flank ham short loin tenderloin shoulder tail cupim porchetta boudin pork chop cow. Ham brisket pastrami salami short loin, frankfurter ribeye.'
$ws.Range("C15").Value = 'local road maintenance'
$ws.Range("E15").Value = 'Pastrami tri-tip frankfurter venison salami, kielbasa turducken shoulder andouille jowl pig beef ribs.'
$ws.Range("C16").Value = 'the mayor'
$ws.Range("E16").Value = 'Swine meatloaf cow sirloin, alcatra pork andouille ham t-bone prosciutto cupim turkey jowl turducken tri-tip. Salami sirloin swine chuck landjaeger jerky kielbasa hamburger jowl ground round shank cupim. '
$ws.Range("C17").Value = 'maintenance of inter-municipal roads and highways'
$ws.Range("E17").Value = 'Pork belly tenderloin pork loin biltong, fatback short ribs turducken. Ham hock filet mignon drumstick andouille, turkey tenderloin corned beef landjaeger strip steak jerky.'
$ws.Range("C18").Value = 'nature and species conservation'
$ws.Range("E18").Value = 'Short ribs ribeye shank salami turkey, jerky beef ribs shoulder. Ribeye fatback ham jowl picanha meatball pig drumstick biltong andouille turkey.'
$ws.Range("C19").Value = 'municipal parks and squares'
$ws.Range("E19").Value = 'Tri-tip biltong flank, ball tip beef tongue tail pork chop strip steak chuck venison meatloaf jerky corned beef. '
$ws.Range("C20").Value = 'phone and postal services'
$ws.Range("E20").Value = 'Doner prosciutto meatloaf beef. Pork chop tri-tip doner shank chuck alcatra brisket boudin jerky ribeye shoulder biltong porchetta hamburger.'
$ws.Range("C21").Value = 'kosovo police'
$ws.Range("E21").Value = 'T-bone chicken short ribs hamburger, salami pig jerky leberkas ham prosciutto. Kevin frankfurter kielbasa turducken turkey strip steak beef flank meatloaf landjaeger doner.'
$ws.Range("C22").Value = 'access to and quality of preschool education'
$ws.Range("E22").Value = 'Capicola cupim hamburger ham hock ribeye t-bone shoulder pork chop cow short loin tenderloin biltong prosciutto ground round kevin.'
$ws.Range("C23").Value = 'primary and secondary schools/education'
$ws.Range("E23").Value = 'Sausage brisket pork, pancetta leberkas fatback doner shankle pastrami swine tongue salami ground round beef ribs.'
$ws.Range("C24").Value = 'public procurement/tenders'
$ws.Range("E24").Value = 'Turkey alcatra ham shank pork chop. Sirloin spare ribs beef ribs capicola, leberkas salami meatball shoulder turkey swine frankfurter.'
$ws.Range("C25").Value = 'recruitment of local staff'
$ws.Range("E25").Value = 'Corsair American Main yo-ho-ho league ho lateen sail splice the main brace skysail pressgang fire ship. Sutler carouser galleon gally American Main Jack Ketch ye chantey yard no prey, no pay. '
$ws.Range("C26").Value = 'sewage and sanitation'
$ws.Range("E26").Value = 'Gibbet killick cable Admiral of the Black league starboard loaded to the gunwalls Sail ho ahoy yawl.'
$ws.Range("C27").Value = 'sidewalks'
$ws.Range("E27").Value = 'Measured fer yer chains list chase guns maroon Davy Jones'' Locker Cat o''nine tails heave to long boat rope''s end execution dock. '
$ws.Range("C28").Value = 'social services'
$ws.Range("E28").Value = 'American Main loaded to the gunwalls fire ship league clipper long boat topmast Sink me tackle cutlass. Loot come about lee coxswain gunwalls stern aye Buccaneer bilge water topmast.'
$ws.Range("C29").Value = 'sports activities'
$ws.Range("E29").Value = 'Red ensign Spanish Main hornswaggle warp trysail gun galleon tender scuppers hulk. Arr gun jolly boat lugsail rope''s end capstan black spot clap of thunder port nipper. '
$ws.Range("C30").Value = 'procedures for tax payment'
$ws.Range("E30").Value = 'Ho heave down scuttle walk the plank scallywag rutters grog blossom hang the jib Jack Tar pirate.'
$ws.Range("C31").Value = 'traffic and parking control and regulation'
$ws.Range("E31").Value = 'You know why the yankees always win, frank? it''s ''cause the other teams can''t stop staring at those damn pinstripes.'
$ws.Range("C32").Value = 'public transport'
$ws.Range("E32").Value = 'Two little mice fell in a bucket of cream. the first mouse quickly gave up and drowned. the second mouse, wouldn''t quit. he struggled so hard that eventually he churned that cream into butter and crawled out.'
$ws.Range("C33").Value = 'urban and rural planning'
$ws.Range("E33").Value = 'Croque monsieur macaroni cheese cheesy feet. Cheesy grin fondue croque monsieur emmental cheese and wine halloumi croque monsieur cheese and wine. '
$ws.Range("C34").Value = 'waste collection services'
$ws.Range("E34").Value = 'Fromage roquefort cheeseburger. Stinking bishop cheddar pepper jack hard cheese camembert de normandie queso edam pecorino. '
$ws.Range("C35").Value = 'water supply'
$ws.Range("E35").Value = 'Çhëësë slices fromage cheese on toast st. agur blue cheese cheese strings cow pepper jack lancashire. Stilton cheese triangles.'
$ws.Range("C36").Value = 'youth activities'
$ws.Range("E36").Value = 'Lo-fi kombuçha PBR&B four loko williamsburg strëët art messenger bag, everyday carry literally put a bird on it meggings vice stumptown. Ethical kombucha affogato, tattooed plaid photo booth readymade. '

$ws = $wb.Worksheets.Item("problems")
$ws.Range("C2").Value = 'Çorruption'
$ws.Range("E2").Value = 'Chalk and cheese bocconcini macaroni cheese rubber cheese roquefort the big cheese cheesecake dolcelatte. Cheesecake pepper jack bocconcini squirty cheese.'
$ws.Range("C3").Value = 'Crimë'
$ws.Range("E3").Value = 'Gouda gouda cheese strings. Jarlsberg cheese and biscuits mozzarella squirty cheese brie cheddar manchego feta.'
$ws.Range("C4").Value = 'Environmental pollution'
$ws.Range("E4").Value = 'No wait, Doc, the bruise, the bruise on your head, I know how that happened, you told me the whole story. you were standing on your toilet and you were hanging a clock, and you fell, and you hit your head on the sink, and that''s when you came up with the idea for the flux capacitor, which makes time travel possible.'
$ws.Range("C5").Value = 'Inter-ethnic relations'
$ws.Range("E5").Value = 'About 30 years, it''s a nice round number. Don''t worry, I''ll take care of the lightning, you take care of your pop.'
$ws.Range("C6").Value = 'Lack of eçonomic growth'
$ws.Range("E6").Value = 'By the way, what happened today, did he ask her out? Get out of town, I didn''t know you did anything creative. Ah, let me read some. Who''s are these?'
$ws.Range("C7").Value = 'Lack of general or personal security'
$ws.Range("E7").Value = 'Yeah. Well, Marty, I want to thank you for all your good advise, I''ll never forget it. What''s the meaning of this. Look, you gotta listen to me. Yeah but George, Lorraine wants to go with you. Give her a break.'
$ws.Range("C8").Value = 'Limited freedom of movement'
$ws.Range("E8").Value = 'He''s your brother, Mom. The storm. Yeah, well, I still don''t understand what Dad was doing in the middle of the street. Maybe you were adopted. Alright, okay Jennifer.'
$ws.Range("C9").Value = 'Poor electricity supply'
$ws.Range("E9").Value = 'What if I send in the tape and they don''t like it. I mean, what if they say I''m no good. What if they say, ''Get out of here, kid, you got no future.'' I mean, I just don''t think I can take that kind of rejection. Jesus, I''m beginning to sound like my old man.'
$ws.Range("C10").Value = 'Poor water supply'
$ws.Range("E10").Value = 'Why is she gonna get angry with you? Calvin, why do you keep calling me Calvin? That''s right. What were you doing in the middle of the street, a kid your age.'
$ws.Range("C11").Value = 'Poorly functioning rubbish collection service '
$ws.Range("E11").Value = 'Now remember, according to my theory you interfered with with your parent''s first meeting.'
$ws.Range("C12").Value = 'Poverty/Low standard of living'
$ws.Range("E12").Value = 'They don''t mëët, they don''t fall in love, they won''t get married and they wont have kids. '
$ws.Range("C13").Value = 'Road infrastructurë'
$ws.Range("E13").Value = 'That''s why your older brother''s disappeared from that photograph. Your sister will follow and unless you repair the damages, you will be next.'
$ws.Range("C14").Value = 'Unemployment'
$ws.Range("E14").Value = 'Manchëgo rëd leicëster çaerphilly. Cow stinking bishop fromage frais dolcelatte red leicester cheesecake cheesy feet babybel. Caerphilly port-salut cheesy feet fondue when the cheese comes out everybody''s happy goat paneer cheddar. '

$ws_i = $wb.Worksheets.Item("indicators")
$ws_p = $wb.Worksheets.Item("problems")
$ws_i.Range("E1").Select()
$ws_p.Activate()
$ws_p.Range("E14").Select()